$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of data (row 55)
$ws.Cells.Item(55, 1).Value = 1948
$ws.Cells.Item(55, 2).Value = "Delete Duplicate Folders in System"
$ws.Cells.Item(55, 3).Value = "Trie"
$ws.Cells.Item(55, 4).Value = "Make a private class: name, subfolders, serial -> Build a Trie -> Post Order to Serialize -> Postorder to Mark Delete -> Preorder to collect"

# Update the view so the new row is visible/selected
$ws.Application.ActiveWindow.ScrollRow = 31
$ws.Range("D55").Select()
